# "Generate Report for Handoff"
#
# Appends two new handoff-report rows (a second .png dependency pair plus
# the .md "Include" source) to each of the three worksheets:
#   Overview (per-file summary), zh-cn (detail table), de-de (detail table)
#
# Row 2 on every sheet already describes 5484ffd3-...png (re-stamped with a
# newer handoff timestamp); rows 3/4 are brand new for 5ab3957f-...md and
# bdbe82f8-...png.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newHandoffDate = "2016-03-24 21:19:17"
$zhHandoffDatetime = "2016-03-24 21:19:12"
$epoch = "0001-01-01 00:00:00"

$pngFile1 = "5484ffd3-cfc6-4682-8f43-fe02410a9677.png"
$mdFile = "5ab3957f-b06b-4999-bc84-69007e006f04.md"
$pngFile2 = "bdbe82f8-65f0-4264-b244-1551f3ff0d78.png"

$zhPngTarget1 = "6e2359d2c2fd07690baa5435af5a7ab9e915f483.png"
$zhMdXlf = "5ab3957f-b06b-4999-bc84-69007e006f04.9b4ec142bcabd59b8b47a92b311877e6f80a3805.zh-cn.xlf"
$zhPngTarget2 = "441739558ef6de789c4737b273aed9e7e74e0682.png"

$deMdXlf = "5ab3957f-b06b-4999-bc84-69007e006f04.9b4ec142bcabd59b8b47a92b311877e6f80a3805.de-de.xlf"

$dependencyFrom = "e2e\5ab3957f-b06b-4999-bc84-69007e006f04.md"

$dateTimeFormat = "yyyy-mm-dd HH:mm:ss"

$e2eBase = "https://github.com/OpenLocalizationTest/oltest/blob/8c904b3ad8031658e672d27c4eb83524f1333844/e2e/"
$zhHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0cfb4c7462b1edddcbbcdace85eaebe394e808b9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e754c27d8d26fb94fb52fc38bc0f48360a20e3b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------

# NOTE: Hyperlinks.Delete() wipes every hyperlink on the sheet (not just
# the scoped range), so it must be called exactly once, up front, before
# any of the fresh Hyperlinks.Add() calls below.
$ws1.Hyperlinks.Delete()

# Row 2 keeps the same source file but the handoff got re-run -> new date.
$ws1.Hyperlinks.Add($ws1.Range("A2"), ($e2eBase + $pngFile1), [Type]::Missing, [Type]::Missing, $pngFile1) | Out-Null
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = $newHandoffDate

$ws1.Hyperlinks.Add($ws1.Range("A3"), ($e2eBase + $mdFile), [Type]::Missing, [Type]::Missing, $mdFile) | Out-Null
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = $newHandoffDate
$ws1.Range("D3").NumberFormat = $dateTimeFormat

$ws1.Hyperlinks.Add($ws1.Range("A4"), ($e2eBase + $pngFile2), [Type]::Missing, [Type]::Missing, $pngFile2) | Out-Null
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = $newHandoffDate
$ws1.Range("D4").NumberFormat = $dateTimeFormat

# ---------------------------------------------------------------------
# zh-cn / de-de detail sheets:
# A Source File Name | B File Extension | C Status | D Latest Handoff File
# E Latest Handoff Datetime | H Latest Handback DateTime
# J Handoff Reason | K Dependency From
# ---------------------------------------------------------------------

function Fill-DetailSheet {
    param(
        $ws,
        $htBase,
        $handoffDatetime,
        $zhTarget1,
        $mdXlf,
        $zhTarget2
    )

    # Row 2: 5484ffd3-...png - depends on the .md handed off alongside it.
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), ($e2eBase + $pngFile1), [Type]::Missing, [Type]::Missing, $pngFile1) | Out-Null
    $ws.Range("B2").Value = ".png"
    $ws.Range("C2").Value = "Ready for handoff"
    $ws.Range("D2").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("D2"), ($htBase + $zhTarget1), [Type]::Missing, [Type]::Missing, $zhTarget1) | Out-Null
    $ws.Range("E2").Value = $handoffDatetime
    $ws.Range("E2").NumberFormat = $dateTimeFormat
    $ws.Range("H2").Value = $epoch
    $ws.Range("H2").NumberFormat = $dateTimeFormat
    $ws.Range("J2").Value = "IsDependency"
    $ws.Range("K2").Value = $dependencyFrom

    # Row 3: 5ab3957f-...md - the source file itself, directly included.
    $ws.Hyperlinks.Add($ws.Range("A3"), ($e2eBase + $mdFile), [Type]::Missing, [Type]::Missing, $mdFile) | Out-Null
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D3"), ($htBase + $mdXlf), [Type]::Missing, [Type]::Missing, $mdXlf) | Out-Null
    $ws.Range("E3").Value = $handoffDatetime
    $ws.Range("E3").NumberFormat = $dateTimeFormat
    $ws.Range("H3").Value = $epoch
    $ws.Range("H3").NumberFormat = $dateTimeFormat
    $ws.Range("J3").Value = "Include"

    # Row 4: bdbe82f8-...png - another dependency of the same .md file.
    $ws.Hyperlinks.Add($ws.Range("A4"), ($e2eBase + $pngFile2), [Type]::Missing, [Type]::Missing, $pngFile2) | Out-Null
    $ws.Range("B4").Value = ".png"
    $ws.Range("C4").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D4"), ($htBase + $zhTarget2), [Type]::Missing, [Type]::Missing, $zhTarget2) | Out-Null
    $ws.Range("E4").Value = $handoffDatetime
    $ws.Range("E4").NumberFormat = $dateTimeFormat
    $ws.Range("H4").Value = $epoch
    $ws.Range("H4").NumberFormat = $dateTimeFormat
    $ws.Range("J4").Value = "IsDependency"
    $ws.Range("K4").Value = $dependencyFrom
}

Fill-DetailSheet $ws2 $zhHtBase $zhHandoffDatetime $zhPngTarget1 $zhMdXlf $zhPngTarget2
Fill-DetailSheet $ws3 $deHtBase $newHandoffDate $zhPngTarget1 $deMdXlf $zhPngTarget2

Write-Host "Report rows appended to Overview, zh-cn, de-de"
